# "start prefix wrangle, add data progress function"
#
# Restructure the 地词(place-name) merge sheet:
#   - row 7:  add a lone space marker in F7
#   - row 10: add "素浐" as an alias in D10
#   - row 12: prefix a new "长安" alias in A12 (pushing the existing
#             帝城/帝京/... group one column right) and fold the
#             old 京/京城/... alias group (previously its own row 16)
#             onto the end of row 12
#   - the old row 16 is then removed, and everything below it shifts
#     up by one row
#   - update the saved selection to C18
#
# NOTE: this engine's Range.Value getter (a COM ParameterizedProperty)
# does not evaluate correctly here -- it yields the reflected property
# descriptor string instead of the cell's contents. Range.Value2 works
# correctly for both reads and writes, so it is used throughout instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 7: new trailing blank-ish alias (brand-new cell, default style) --
$ws.Range("F7").Value2 = " "

# --- row 10: new alias "素浐" (matches the row's existing s="1" style) -----
$ws.Range("A10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value2 = "素浐"

# --- row 12: remember the 京-group currently living on row 16 so it can
#     be appended after we shift row 12's own contents one column right
$jingCheng = $ws.Range("A16").Value2
$jing      = $ws.Range("B16").Value2
$jingChi   = $ws.Range("C16").Value2
$jingDian  = $ws.Range("D16").Value2
$jingDu    = $ws.Range("E16").Value2
$jingGuo   = $ws.Range("F16").Value2
$jingHua   = $ws.Range("G16").Value2
$jingYi    = $ws.Range("H16").Value2
$jingLuo   = $ws.Range("I16").Value2

# remember row 12's own current A:E contents before overwriting them
$diCheng = $ws.Range("A12").Value2
$diJing  = $ws.Range("B12").Value2
$diJu    = $ws.Range("C12").Value2
$diLi    = $ws.Range("D12").Value2
$diXiang = $ws.Range("E12").Value2

# build the whole new row 12 (A:O), copying A12's format (s="1") onto
# every cell first so newly-touched columns (F:O) pick up the same style
$ws.Range("A12").Copy()
$ws.Range("A12:O12").PasteSpecial(-4122)

$ws.Range("A12").Value2 = "长安"
$ws.Range("B12").Value2 = $diCheng
$ws.Range("C12").Value2 = $diJing
$ws.Range("D12").Value2 = $diJu
$ws.Range("E12").Value2 = $diLi
$ws.Range("F12").Value2 = $diXiang
$ws.Range("G12").Value2 = $jingCheng
$ws.Range("H12").Value2 = $jing
$ws.Range("I12").Value2 = $jingChi
$ws.Range("J12").Value2 = $jingDian
$ws.Range("K12").Value2 = $jingDu
$ws.Range("L12").Value2 = $jingGuo
$ws.Range("M12").Value2 = $jingHua
$ws.Range("N12").Value2 = $jingYi
$ws.Range("O12").Value2 = $jingLuo

$excel.CutCopyMode = 0

# --- remove the now-redundant old row 16 (its data now lives on row 12);
#     everything below shifts up by one row to close the gap
$ws.Rows.Item(16).Delete()

# --- restore the saved selection -------------------------------------
$ws.Range("C18").Select()
